$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date/time number formatting (styles) from row 18 so the new
# rows reuse the same cellXfs entries (s="7" for C/D, s="6" for F),
# without touching any other columns.
$ws.Range("C18").Copy($ws.Range("C19:C24"))
$ws.Range("D18").Copy($ws.Range("D19:D24"))
$ws.Range("F18").Copy($ws.Range("F19:F24"))

function Set-PlainCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value2 = $value
    $cell.Style = "Normal"
}

# Row 19
Set-PlainCell 19 1 18
Set-PlainCell 19 2 "Snakes"
$ws.Cells.Item(19, 3).Value2 = 45408.91874230324
$ws.Cells.Item(19, 4).Value2 = 45408.91880518518
$ws.Cells.Item(19, 6).Value2 = 0.00004629629629629629
Set-PlainCell 19 9 "Team1"
Set-PlainCell 19 10 "Process1"
Set-PlainCell 19 11 "Person1"

# Row 20
Set-PlainCell 20 1 19
Set-PlainCell 20 2 "Snakes"
$ws.Cells.Item(20, 3).Value2 = 45408.91925273148
$ws.Cells.Item(20, 4).Value2 = 45408.91981146991
$ws.Cells.Item(20, 6).Value2 = 0.0005439814814814814
Set-PlainCell 20 9 "Team1"
Set-PlainCell 20 10 "Process1"
Set-PlainCell 20 11 "Person1"

# Row 21
Set-PlainCell 21 1 20
Set-PlainCell 21 2 "Snakes"
$ws.Cells.Item(21, 3).Value2 = 45408.93170967593
$ws.Cells.Item(21, 4).Value2 = 45408.93189614583
$ws.Cells.Item(21, 6).Value2 = 0.0001851851851851852
Set-PlainCell 21 7 "CTkCheckBox"
Set-PlainCell 21 9 "Team1"
Set-PlainCell 21 10 "Process1"
Set-PlainCell 21 11 "Person1"

# Row 22
Set-PlainCell 22 1 21
Set-PlainCell 22 2 "Snakes"
$ws.Cells.Item(22, 3).Value2 = 45408.95428038194
$ws.Cells.Item(22, 4).Value2 = 45408.9543390625
$ws.Cells.Item(22, 6).Value2 = 0.00005787037037037037
Set-PlainCell 22 7 "Measure other"
Set-PlainCell 22 9 "Team1"
Set-PlainCell 22 10 "Process1"
Set-PlainCell 22 11 "Person1"

# Row 23
Set-PlainCell 23 1 22
Set-PlainCell 23 2 "Snakes"
$ws.Cells.Item(23, 3).Value2 = 45408.95437648148
$ws.Cells.Item(23, 4).Value2 = 45408.9544322338
$ws.Cells.Item(23, 6).Value2 = 0.00003472222222222222
Set-PlainCell 23 7 "Measure other"
Set-PlainCell 23 9 "Team1"
Set-PlainCell 23 10 "Process1"
Set-PlainCell 23 11 "Person1"

# Row 24
Set-PlainCell 24 1 23
Set-PlainCell 24 2 "Snakes"
$ws.Cells.Item(24, 3).Value2 = 45410.65324151301
$ws.Cells.Item(24, 4).Value2 = 45410.65335950792
$ws.Cells.Item(24, 6).Value2 = 0.0001157407407407407
Set-PlainCell 24 9 "Team1"
Set-PlainCell 24 10 "Process1"
Set-PlainCell 24 11 "Person1"

# Update C18/D18 precision as in the diff (tiny recompute drift)
$ws.Cells.Item(18, 3).Value2 = 45408.91043981481
$ws.Cells.Item(18, 4).Value2 = 45408.91049984954

Write-Host $ws.UsedRange.Address()
